$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply every cell update from the source diff in row order. For "Price" (D)
# values that are plain decimal numbers, Excel would otherwise auto-convert the
# assigned text into a numeric value (dropping formatting such as trailing
# zeros). Force those specific cells to remain text by switching their number
# format to Text ("@") before the assignment, then restore the default "Normal"
# style afterwards so no stray formatting remains on the cell.

$ws.Range("D2").Value = '43.765.41'
$ws.Range("E2").Value = '  -1.01%  '
$ws.Range("D3").Value = '2.346.07'
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.673'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '239.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.65%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.31'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.36%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.599'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.100'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.58'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.77'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.33'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.108'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").Value = '2.695.53'
$ws.Range("E15").Value = '  -0.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.56%  '
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("D18").Value = '2.350.03'
$ws.Range("D19").Value = '43.683.81'
$ws.Range("E19").Value = '  -1.43%  '
$ws.Range("E20").Value = '  -1.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.74'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '77.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.53'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("E24").Value = '  +22.52%  '
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.63'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.28'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.63'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '177.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.131'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.32%  '
$ws.Range("E33").Value = '  +3.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0760'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.45%  '
$ws.Range("E35").Value = '  -3.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.50'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.18%  '
$ws.Range("E37").Value = '  -1.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.28'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0281'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '68.95'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +31.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.112'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +11.54%  '
$ws.Range("E43").Value = '  +1.87%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.13'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.22%  '
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.203'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.82'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.67%  '
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("E49").Value = '  -1.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '99.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.07%  '
